# "updated results with new values"
# Refresh the quicksort benchmark results on the worksheet with the newly
# measured timings (ns) for each input-size column (2^8 .. 2^16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quicksort")

$bestCase          = @(6900, 13720, 28380, 69960, 154600, 432150, 1347540, 4532630, 16446350)
$averageCase       = @(1720, 3170, 6120, 12950, 29800, 59470, 125160, 272120, 566800)
$worstCaseSorted   = @(7060, 25700, 97790, 382030, 1520610, 6083490, 24107830, 96083910, 385193780)
$worstCaseAllSame  = @(14660, 51990, 195320, 767810, 3032800, 12073520, 48234490, 192396910, 772075630)

for ($i = 0; $i -lt 9; $i++) {
    $col = $i + 2
    $ws.Cells.Item(2, $col).Value = $bestCase[$i]
    $ws.Cells.Item(3, $col).Value = $averageCase[$i]
    $ws.Cells.Item(4, $col).Value = $worstCaseSorted[$i]
    $ws.Cells.Item(5, $col).Value = $worstCaseAllSame[$i]
}

$null = $ws.Range("A1:J3").Select()

$wb.Save()
